$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.490.01'
$ws.Range("E2").Value = '  +1.89%  '
$ws.Range("D3").Value = '1.859.17'
$ws.Range("E3").Value = '  +0.74%  '
$ws.Range("E4").Value = '  -0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.16'
$ws.Range("E5").Value = '  +0.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.009'
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4762'
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3784'
$ws.Range("E8").Value = '  +2.85%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07318'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9289'
$ws.Range("E10").Value = '  -0.35%  '
$ws.Range("E11").Value = '  +4.13%  '
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("D13").Value = '1.870.04'
$ws.Range("E13").Value = '  +1.20%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.435'
$ws.Range("E14").Value = '  +0.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.554'
$ws.Range("E15").Value = '  +1.30%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.06'
$ws.Range("E16").Value = '  +1.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.011'
$ws.Range("E17").Value = '  -0.39%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008809'
$ws.Range("E18").Value = '  +1.70%  '
$ws.Range("E19").Value = '  -0.41%  '
$ws.Range("D20").Value = '27.513.15'
$ws.Range("E20").Value = '  +1.86%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.65'
$ws.Range("E21").Value = '  +0.74%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.082'
$ws.Range("E22").Value = '  +0.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.70'
$ws.Range("E23").Value = '  +0.46%  '
$ws.Range("E24").Value = '  -1.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '155.44'
$ws.Range("E25").Value = '  +1.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.45'
$ws.Range("E26").Value = '  +1.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.002'
$ws.Range("E27").Value = '  -0.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.943'
$ws.Range("E29").Value = '  -0.54%  '
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("E31").Value = '  +1.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.203'
$ws.Range("E32").Value = '  +1.99%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7521'
$ws.Range("E33").Value = '  +1.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.582'
$ws.Range("E34").Value = '  +1.60%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.706'
$ws.Range("E35").Value = '  +1.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02041'
$ws.Range("E36").Value = '  +3.92%  '
$ws.Range("E37").Value = '  +0.55%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5558'
$ws.Range("E38").Value = '  +5.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05319'
$ws.Range("E39").Value = '  +0.98%  '
$ws.Range("E40").Value = '  +0.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.027'
$ws.Range("E41").Value = '  +0.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.466'
$ws.Range("E42").Value = '  +1.96%  '
$ws.Range("E43").Value = '  +0.52%  '
$ws.Range("E44").Value = '  +0.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4862'
$ws.Range("E45").Value = '  +2.66%  '
$ws.Range("E46").Value = '  -0.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '103.94'
$ws.Range("E47").Value = '  +2.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.665'
$ws.Range("E48").Value = '  +3.44%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '67.28'
$ws.Range("E49").Value = '  +2.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9108'
$ws.Range("E51").Value = '  +2.11%  '
